# Update the "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Row -> new value mapping (only rows that changed)
$updates = @{
    2 = 2158
    3 = 1644
    5 = 1050
    6 = 592
    7 = 31
    8 = 5722
    9 = 82
}

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
